# Saldo_guide.xlsx update:
#  - rename the sheet to reflect the new export timestamp
#  - bump every "Dt. Referencia" (column G) date forward by one day
#    (2024-06-05 -> 2024-06-06)
#  - correct the balance for the account in row 245 (Vl. Projetado /
#    Vl. Total): 4169.91 -> 14169.91

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet to match the new export run.
$ws.Name = "IClientBalance-20240606-100000-"

# Data rows are 2..257 inclusive; column G (7) holds the reference date.
for ($r = 2; $r -le 257; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $cur = $cell.Value()
    $cell.Value = $cur.AddDays(1)
}

# Row 245 balance correction (D = Vl. Projetado, H = Vl. Total).
$ws.Cells.Item(245, 4).Value = 14169.91
$ws.Cells.Item(245, 8).Value = 14169.91
